$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; this shifts the existing rows 50-116
# down to 51-117 (and grows the used range from A1:T116 to A1:T117).
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new price-record data.
$ws.Cells.Item(50, 1).Value = 1
$ws.Cells.Item(50, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value = 44790
$ws.Cells.Item(50, 5).Value = 15
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100106
$ws.Cells.Item(50, 8).Value = "Oleaginosos"
$ws.Cells.Item(50, 9).Value = 100106002
$ws.Cells.Item(50, 10).Value = "Palta"
$ws.Cells.Item(50, 11).Value = "Hass"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 400
$ws.Cells.Item(50, 14).Value = 24000
$ws.Cells.Item(50, 15).Value = 25000
$ws.Cells.Item(50, 16).Value = 24500
$ws.Cells.Item(50, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(50, 18).Value = "Perú"
$ws.Cells.Item(50, 19).Value = 2450
$ws.Cells.Item(50, 20).Value = 10
